$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.8494554131357859
$ws.Range("C2").Value = 0.2321382534920531
$ws.Range("D2").Value = 0.2214068020382598
$ws.Range("F2").Value = 1.501654300016142
$ws.Range("G2").Value = 0.8614496656753516
$ws.Range("H2").Value = 0.8890446733713873
$ws.Range("I2").Value = 0.6420947543634981
$ws.Range("J2").Value = 0.2547594434794576
$ws.Range("M2").Value = 0.4393645854472794
$ws.Range("N2").Value = 1.25749254344748
$ws.Range("B3").Value = 0.7631781101640343
$ws.Range("C3").Value = 0.205584274731109
$ws.Range("D3").Value = 0.2201062118342776
$ws.Range("F3").Value = 1.49282660655598
$ws.Range("G3").Value = 0.8508640101280349
$ws.Range("H3").Value = 0.8895383111857882
$ws.Range("I3").Value = 0.6450442097021565
$ws.Range("J3").Value = 0.2548074062483963
$ws.Range("M3").Value = 0.4145007313208211
$ws.Range("N3").Value = 1.273565970302339
$ws.Range("B4").Value = 0.7103582971338653
$ws.Range("C4").Value = 0.1893067507571686
$ws.Range("D4").Value = 0.2193745584009861
$ws.Range("F4").Value = 1.488372364086359
$ws.Range("G4").Value = 0.845040689415768
$ws.Range("H4").Value = 0.89035154534659
$ws.Range("I4").Value = 0.6472900650962927
$ws.Range("J4").Value = 0.2549627622075974
$ws.Range("M4").Value = 0.39941871973069
$ws.Range("N4").Value = 1.283958416442207
$ws.Range("B5").Value = 0.6888733468423993
$ws.Range("C5").Value = 0.1826802941865253
$ws.Range("D5").Value = 0.2190932877501126
$ws.Range("F5").Value = 1.486799738988438
$ws.Range("G5").Value = 0.8428371493653373
$ws.Range("H5").Value = 0.8908110888470731
$ws.Range("I5").Value = 0.6483143871300321
$ws.Range("J5").Value = 0.2550577152243747
$ws.Range("M5").Value = 0.3933192205790874
$ws.Range("N5").Value = 1.288324933042302
$ws.Range("B6").Value = 0.6853081947205339
$ws.Range("C6").Value = 0.1815803862368739
$ws.Range("D6").Value = 0.2190476044403411
$ws.Range("F6").Value = 1.486553241675864
$ws.Range("G6").Value = 0.8424814750574683
$ws.Range("H6").Value = 0.8908951313812707
$ws.Range("I6").Value = 0.6484910585826036
$ws.Range("J6").Value = 0.2550753928745948
$ws.Range("M6").Value = 0.392309218965444
$ws.Range("N6").Value = 1.289057930533561
$ws.Range("B7").Value = 0.7100683825997578
$ws.Range("C7").Value = 0.1892173565840665
$ws.Range("D7").Value = 0.2193706966464433
$ws.Range("F7").Value = 1.488350173705257
$ws.Range("G7").Value = 0.8450102860991819
$ws.Range("H7").Value = 0.8903572242434876
$ws.Range("I7").Value = 0.647303437963199
$ws.Range("J7").Value = 0.2549639146723308
$ws.Range("M7").Value = 0.3993362711003101
$ws.Range("N7").Value = 1.284016772435688
$ws.Range("B8").Value = 0.8196752175433062
$ws.Range("C8").Value = 0.2229769134956427
$ws.Range("D8").Value = 0.220944504388477
$ws.Range("F8").Value = 1.498409774300526
$ws.Range("G8").Value = 0.8576590226554401
$ws.Range("H8").Value = 0.8891089166934449
$ws.Range("I8").Value = 0.6430213244354661
$ws.Range("J8").Value = 0.2547498387192846
$ws.Range("M8").Value = 0.4307533175213365
$ws.Range("N8").Value = 1.262926005174361
$ws.Range("B9").Value = 1.035825911909797
$ws.Range("C9").Value = 0.2893928073038694
$ws.Range("D9").Value = 0.2245594539516134
$ws.Range("F9").Value = 1.525822434980782
$ws.Range("G9").Value = 0.887857070102342
$ws.Range("H9").Value = 0.890715963792104
$ws.Range("I9").Value = 0.6380858909595943
$ws.Range("J9").Value = 0.2553300894400294
$ws.Range("M9").Value = 0.4938229642146439
$ws.Range("N9").Value = 1.225722161012314
$ws.Range("B10").Value = 1.195365574067637
$ws.Range("C10").Value = 0.3383258001064178
$ws.Range("D10").Value = 0.2275351832131349
$ws.Range("F10").Value = 1.550681970061206
$ws.Range("G10").Value = 0.9133739650067128
$ws.Range("H10").Value = 0.8943804402793774
$ws.Range("I10").Value = 0.6365865838166656
$ws.Range("J10").Value = 0.2563679577723406
$ws.Range("M10").Value = 0.5410522133683884
$ws.Range("N10").Value = 1.200924912900362
$ws.Range("B11").Value = 1.268103655378525
$ws.Range("C11").Value = 0.3606183352541734
$ws.Range("D11").Value = 0.2289578744912717
$ws.Range("F11").Value = 1.56302366878235
$ws.Range("G11").Value = 0.9257147510579102
$ws.Range("H11").Value = 0.8965895110360407
$ws.Range("I11").Value = 0.6363699144211026
$ws.Range("J11").Value = 0.2569733526864368
$ws.Range("M11").Value = 0.5627323655541971
$ws.Range("N11").Value = 1.190195547516073
$ws.Range("B12").Value = 1.295670699040443
$ws.Range("C12").Value = 0.3690646756544993
$ws.Range("D12").Value = 0.2295064851542605
$ws.Range("F12").Value = 1.567846222949484
$ws.Range("G12").Value = 0.9304939744981766
$ws.Range("H12").Value = 0.8975041705814988
$ws.Range("I12").Value = 0.6363550693918327
$ws.Range("J12").Value = 0.2572217924542031
$ws.Range("M12").Value = 0.5709701013507242
$ws.Range("N12").Value = 1.186211958902735
$ws.Range("B13").Value = 1.289732643697789
$ws.Range("C13").Value = 0.3672453995968681
$ws.Range("D13").Value = 0.229387894025848
$ws.Range("F13").Value = 1.5668009641971
$ws.Range("G13").Value = 0.9294599572364746
$ws.Range("H13").Value = 0.8973037044094667
$ws.Range("I13").Value = 0.6363552731451563
$ws.Range("J13").Value = 0.2571674326003048
$ws.Range("M13").Value = 0.569194717696206
$ws.Range("N13").Value = 1.187066362650647
$ws.Range("B14").Value = 1.270371161385981
$ws.Range("C14").Value = 0.3613131280659445
$ws.Range("D14").Value = 0.2290028115650244
$ws.Range("F14").Value = 1.563417433854781
$ws.Range("G14").Value = 0.9261058119315919
$ws.Range("H14").Value = 0.8966631936907561
$ws.Range("I14").Value = 0.6363673446392397
$ws.Range("J14").Value = 0.2569934072278102
$ws.Range("M14").Value = 0.5634095301071085
$ws.Range("N14").Value = 1.18986622331585
$ws.Range("B15").Value = 1.258514631858759
$ws.Range("C15").Value = 0.3576800432186928
$ws.Range("D15").Value = 0.228768220927023
$ws.Range("F15").Value = 1.561364345037646
$ws.Range("G15").Value = 0.9240651287175012
$ws.Range("H15").Value = 0.8962810429482033
$ws.Range("I15").Value = 0.6363834990056603
$ws.Range("J15").Value = 0.2568893115496707
$ws.Range("M15").Value = 0.5598695687435367
$ws.Range("N15").Value = 1.19159156223888
$ws.Range("B16").Value = 1.19061518453907
$ws.Range("C16").Value = 0.3368695830299941
$ws.Range("D16").Value = 0.2274435909777992
$ws.Range("F16").Value = 1.549896238068172
$ws.Range("G16").Value = 0.9125822695250463
$ws.Range("H16").Value = 0.8942469958711854
$ws.Range("I16").Value = 0.6366101323192694
$ws.Range("J16").Value = 0.2563310776658412
$ws.Range("M16").Value = 0.539639281633626
$ws.Range("N16").Value = 1.201637203125262
$ws.Range("B17").Value = 1.149002253166998
$ws.Range("C17").Value = 0.3241113665854414
$ws.Range("D17").Value = 0.2266486102154772
$ws.Range("F17").Value = 1.543125835687732
$ws.Range("G17").Value = 0.9057260980696213
$ws.Range("H17").Value = 0.893138142742032
$ws.Range("I17").Value = 0.6368685585922904
$ws.Range("J17").Value = 0.2560227700813229
$ws.Range("M17").Value = 0.5272785731395615
$ws.Range("N17").Value = 1.207941116330762
$ws.Range("B18").Value = 1.125082949196553
$ws.Range("C18").Value = 0.3167762564076781
$ws.Range("D18").Value = 0.2261978578617487
$ws.Range("F18").Value = 1.539328878527513
$ws.Range("G18").Value = 0.901851567470473
$ws.Range("H18").Value = 0.8925513720798648
$ws.Range("I18").Value = 0.6370609846511712
$ws.Range("J18").Value = 0.2558579814670239
$ws.Range("M18").Value = 0.5201874150879959
$ws.Range("N18").Value = 1.211618816311553
$ws.Range("B19").Value = 1.11698694177079
$ws.Range("C19").Value = 0.3142932446291127
$ws.Range("D19").Value = 0.2260463588835506
$ws.Range("F19").Value = 1.538059973547149
$ws.Range("G19").Value = 0.9005515424742896
$ws.Range("H19").Value = 0.8923614575350314
$ws.Range("I19").Value = 0.6371336485184145
$ws.Range("J19").Value = 0.2558043401254011
$ws.Range("M19").Value = 0.5177896348660056
$ws.Range("N19").Value = 1.212872924079196
$ws.Range("B20").Value = 1.153430436688723
$ws.Range("C20").Value = 0.3254691823160556
$ws.Range("D20").Value = 0.2267325650042693
$ws.Range("F20").Value = 1.543836493347797
$ws.Range("G20").Value = 0.9064488077996202
$ws.Range("H20").Value = 0.8932509011509353
$ws.Range("I20").Value = 0.6368365149284401
$ws.Range("J20").Value = 0.2560542917338609
$ws.Range("M20").Value = 0.5285924890273535
$ws.Range("N20").Value = 1.207264686032207
$ws.Range("B21").Value = 1.276057487120852
$ws.Range("C21").Value = 0.3630554545006817
$ws.Range("D21").Value = 0.229115652240182
$ws.Range("F21").Value = 1.564407210508008
$ws.Range("G21").Value = 0.9270881229238057
$ws.Range("H21").Value = 0.8968492052668182
$ws.Range("I21").Value = 0.6363619728183991
$ws.Range("J21").Value = 0.2570440016753039
$ws.Range("M21").Value = 0.5651080230815495
$ws.Range("N21").Value = 1.189041680530046
$ws.Range("B22").Value = 1.356333416819552
$ws.Range("C22").Value = 0.3876472831919386
$ws.Range("D22").Value = 0.230730632018961
$ws.Range("F22").Value = 1.578720195453499
$ws.Range("G22").Value = 0.9411954579597932
$ws.Range("H22").Value = 0.8996564060733192
$ws.Range("I22").Value = 0.6364436746487669
$ws.Range("J22").Value = 0.2578027002478649
$ws.Range("M22").Value = 0.5891357777906165
$ws.Range("N22").Value = 1.177594708302085
$ws.Range("B23").Value = 1.31347678287085
$ws.Range("C23").Value = 0.3745197067743788
$ws.Range("D23").Value = 0.2298634446371892
$ws.Range("F23").Value = 1.571001428169737
$ws.Range("G23").Value = 0.9336093207320175
$ws.Range("H23").Value = 0.8981164112804834
$ws.Range("I23").Value = 0.6363641214246485
$ws.Range("J23").Value = 0.2573875238375507
$ws.Range("M23").Value = 0.5762968728316054
$ws.Range("N23").Value = 1.183661773567771
$ws.Range("B24").Value = 1.151428440775021
$ws.Range("C24").Value = 0.3248553145774338
$ws.Range("D24").Value = 0.226694589440342
$ws.Range("F24").Value = 1.543514907767005
$ws.Range("G24").Value = 0.9061218615332081
$ws.Range("H24").Value = 0.8931997650781511
$ws.Range("I24").Value = 0.6368508652843303
$ws.Range("J24").Value = 0.2560400019806295
$ws.Range("M24").Value = 0.5279984203667212
$ws.Range("N24").Value = 1.207570333534086
$ws.Range("B25").Value = 0.9772222197728411
$ws.Range("C25").Value = 0.2714019660237454
$ws.Range("D25").Value = 0.2235251708492143
$ws.Range("F25").Value = 1.517580264919218
$ws.Range("G25").Value = 0.8791058973163075
$ws.Range("H25").Value = 0.889845893924516
$ws.Range("I25").Value = 0.6390487628059915
$ws.Range("J25").Value = 0.2550658592425705
$ws.Range("M25").Value = 0.476604418955823
$ws.Range("N25").Value = 1.235341628552749
